$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers; force Text format first so
# Excel stores the exact literal (e.g. "1.00", "0.650") instead of auto-converting
# them to a number and normalizing the text (e.g. "1", "0.65").
$textCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D42", "D44", "D46", "D47", "D48", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '62.065.67'
$ws.Range('E2').Value = '  +1.54%  '
$ws.Range('D3').Value = '3.001.20'
$ws.Range('E3').Value = '  +1.05%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '544.54'
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').Value = '138.78'
$ws.Range('E6').Value = '  +5.47%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('D8').Value = '2.997.94'
$ws.Range('E8').Value = '  +1.09%  '
$ws.Range('D9').Value = '0.486'
$ws.Range('E9').Value = '  -1.07%  '
$ws.Range('D10').Value = '6.77'
$ws.Range('E10').Value = '  +15.95%  '
$ws.Range('D11').Value = '0.148'
$ws.Range('E11').Value = '  +1.50%  '
$ws.Range('D12').Value = '0.442'
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('E13').Value = '  +0.60%  '
$ws.Range('D14').Value = '33.82'
$ws.Range('E14').Value = '  +0.17%  '
$ws.Range('D15').Value = '3.481.46'
$ws.Range('E15').Value = '  +0.74%  '
$ws.Range('D16').Value = '61.977.00'
$ws.Range('E16').Value = '  +1.32%  '
$ws.Range('D17').Value = '2.995.73'
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('E18').Value = '  -2.08%  '
$ws.Range('D19').Value = '6.53'
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('D20').Value = '463.77'
$ws.Range('E20').Value = '  -0.31%  '
$ws.Range('D21').Value = '13.25'
$ws.Range('E21').Value = '  +1.52%  '
$ws.Range('D22').Value = '0.650'
$ws.Range('E22').Value = '  -2.02%  '
$ws.Range('D23').Value = '7.20'
$ws.Range('E23').Value = '  +4.11%  '
$ws.Range('D24').Value = '79.01'
$ws.Range('E24').Value = '  -0.68%  '
$ws.Range('D25').Value = '12.51'
$ws.Range('E25').Value = '  +5.37%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  +0.46%  '
$ws.Range('D28').Value = '7.58'
$ws.Range('E28').Value = '  -0.47%  '
$ws.Range('E29').Value = '  +6.16%  '
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('D31').Value = '25.30'
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('E32').Value = '  +0.24%  '
$ws.Range('D33').Value = '2.33'
$ws.Range('E33').Value = '  +2.10%  '
$ws.Range('D34').Value = '5.50'
$ws.Range('E34').Value = '  +1.80%  '
$ws.Range('D35').Value = '54.81'
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('D36').Value = '5.81'
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('D37').Value = '449.85'
$ws.Range('E37').Value = '  +1.92%  '
$ws.Range('D38').Value = '0.0803'
$ws.Range('E38').Value = '  +2.33%  '
$ws.Range('D39').Value = '0.0387'
$ws.Range('E39').Value = '  +3.06%  '
$ws.Range('D40').Value = '2.924.07'
$ws.Range('E40').Value = '  -7.18%  '
$ws.Range('E41').Value = '  -1.94%  '
$ws.Range('D42').Value = '8.05'
$ws.Range('E42').Value = '  +0.30%  '
$ws.Range('E43').Value = '  +7.66%  '
$ws.Range('D44').Value = '26.59'
$ws.Range('E44').Value = '  +4.00%  '
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('D46').Value = '0.245'
$ws.Range('E46').Value = '  +2.28%  '
$ws.Range('D47').Value = '0.109'
$ws.Range('E47').Value = '  +1.87%  '
$ws.Range('D48').Value = '1.98'
$ws.Range('E48').Value = '  +2.21%  '
$ws.Range('D49').Value = '114.30'
$ws.Range('E49').Value = '  -1.86%  '
$ws.Range('D50').Value = '0.0₃0497'
$ws.Range('E50').Value = '  +3.09%  '
$ws.Range('E51').Value = '  -3.89%  '

# Restore the default (General) cell style now that the literal text is locked in,
# matching the original workbook formatting for these cells.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
